$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = '''29.611.68'
$ws.Range("E2").Value = '  +3.47%  '
$ws.Range("D3").Value = '''1.606.95'
$ws.Range("E3").Value = '  +2.77%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '''212.62'
$ws.Range("E5").Value = '  +1.18%  '
$ws.Range("D6").Value = '''0.522'
$ws.Range("E6").Value = '  +2.98%  '
$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = '''26.89'
$ws.Range("E8").Value = '  +8.13%  '
$ws.Range("D9").Value = '''43.62'
$ws.Range("E9").Value = '  -1.23%  '
$ws.Range("E10").Value = '  +2.69%  '
$ws.Range("E11").Value = '  +2.66%  '
$ws.Range("D12").Value = '''0.0909'
$ws.Range("E12").Value = '  +1.65%  '
$ws.Range("D13").Value = '''1.838.09'
$ws.Range("E13").Value = '  +2.86%  '
$ws.Range("D14").Value = '''1.620.27'
$ws.Range("E14").Value = '  +3.57%  '
$ws.Range("D15").Value = '''29.627.21'
$ws.Range("E15").Value = '  +3.41%  '
$ws.Range("E16").Value = '  +4.06%  '
$ws.Range("E17").Value = '  +2.40%  '
$ws.Range("D18").Value = '''63.43'
$ws.Range("E18").Value = '  +3.29%  '
$ws.Range("D19").Value = '''240.69'
$ws.Range("E19").Value = '  +6.01%  '
$ws.Range("E20").Value = '  +3.81%  '
$ws.Range("E21").Value = '  +2.06%  '
$ws.Range("D22").Value = '''1.00'
$ws.Range("E22").Value = '  +0.05%  '
$ws.Range("E23").Value = '  +1.90%  '
$ws.Range("D24").Value = '''9.23'
$ws.Range("E24").Value = '  +2.28%  '
$ws.Range("E25").Value = '  +0.67%  '
$ws.Range("D26").Value = '''154.55'
$ws.Range("E26").Value = '  +1.88%  '
$ws.Range("D27").Value = '''15.31'
$ws.Range("E27").Value = '  +3.79%  '
$ws.Range("E28").Value = '  +2.68%  '
$ws.Range("E29").Value = '  +3.13%  '
$ws.Range("E30").Value = '  +0.03%  '
$ws.Range("E32").Value = '  +0.88%  '
$ws.Range("E33").Value = '  +1.46%  '
$ws.Range("E34").Value = '  +4.60%  '
$ws.Range("D35").Value = '''1.412.36'
$ws.Range("E35").Value = '  +0.91%  '
$ws.Range("E36").Value = '  -1.06%  '
$ws.Range("E37").Value = '  +4.92%  '
$ws.Range("D38").Value = '''2.82'
$ws.Range("E38").Value = '  +5.50%  '
$ws.Range("E39").Value = '  +0.31%  '
$ws.Range("E40").Value = '  +2.36%  '
$ws.Range("D41").Value = '''0.540'
$ws.Range("E41").Value = '  +4.55%  '
$ws.Range("E42").Value = '  +2.42%  '
$ws.Range("D43").Value = '''0.0492'
$ws.Range("E43").Value = '  +6.83%  '
$ws.Range("D44").Value = '''54.04'
$ws.Range("E44").Value = '  +27.72%  '
$ws.Range("E45").Value = '  +4.33%  '
$ws.Range("E46").Value = '  +0.02%  '
$ws.Range("D47").Value = '''66.03'
$ws.Range("E47").Value = '  +3.42%  '
$ws.Range("E48").Value = '  +1.34%  '
$ws.Range("D49").Value = '''1.747.03'
$ws.Range("D50").Value = '''0.871'
$ws.Range("E50").Value = '  +0.43%  '
$ws.Range("D51").Value = '''86.61'
$ws.Range("E51").Value = '  +2.10%  '
